# Update price, config and deal data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove discontinued sets (42684, 42688, 42696) ---
$ws.Range("A18:A20").EntireRow.Delete()

# --- Insert 4 new rows for sets 10321, 10339, 10357, 10365 ---
$ws.Range("A2:A5").EntireRow.Insert()
$ws.Range("A2:J5").Style = "Normal"

# --- Insert 1 new row for set 21338 (currently pushes down the 21365 row at position 13) ---
$ws.Range("A13").EntireRow.Insert()
$ws.Range("A13:J13").Style = "Normal"

# --- Populate the brand-new rows with their full data ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("A2").Value = '10321'
$ws.Range("B2").Value = 'Corvette'
$ws.Range("C2").Value = '1210'
$ws.Range("D2").Value = 'N/A'
$ws.Range("E2").Value = 'https://www.lego.com/cdn/cs/set/assets/blt2564f1fe0e59bb78/10321.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F2").Value = 'https://www.lego.com/fr-fr/product/10321'
$ws.Range("G2").Value = ''
$ws.Range("H2").Value = ''
$ws.Range("I2").Value = ''
$ws.Range("J2").Value = ''

$ws.Range("A3").NumberFormat = "@"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("A3").Value = '10339'
$ws.Range("B3").Value = 'Le bureau de poste du Père Noël'
$ws.Range("C3").Value = '1440'
$ws.Range("D3").Value = 'N/A'
$ws.Range("E3").Value = 'https://www.lego.com/cdn/cs/set/assets/blt3c4e5efcccc53a93/10339_Prod.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F3").Value = 'https://www.lego.com/fr-fr/product/10339'
$ws.Range("G3").Value = ''
$ws.Range("H3").Value = ''
$ws.Range("I3").Value = ''
$ws.Range("J3").Value = ''

$ws.Range("A4").NumberFormat = "@"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("A4").Value = '10357'
$ws.Range("B4").Value = 'Shelby Cobra 427 S/C'
$ws.Range("C4").Value = '1241'
$ws.Range("D4").Value = 'N/A'
$ws.Range("E4").Value = 'https://www.lego.com/cdn/cs/set/assets/blt8a933e3230c8710d/10357_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F4").Value = 'https://www.lego.com/fr-fr/product/10357'
$ws.Range("G4").Value = ''
$ws.Range("H4").Value = ''
$ws.Range("I4").Value = ''
$ws.Range("J4").Value = ''

$ws.Range("A5").NumberFormat = "@"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("A5").Value = '10365'
$ws.Range("B5").Value = 'Le bateau pirate du capitaine Jack Sparrow'
$ws.Range("C5").Value = '2862'
$ws.Range("D5").Value = 'N/A'
$ws.Range("E5").Value = 'https://www.lego.com/cdn/cs/set/assets/bltcf20096d15e25f4c/10365_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F5").Value = 'https://www.lego.com/fr-fr/product/10365'
$ws.Range("G5").Value = ''
$ws.Range("H5").Value = ''
$ws.Range("I5").Value = ''
$ws.Range("J5").Value = ''

$ws.Range("A13").NumberFormat = "@"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("A13").Value = '21338'
$ws.Range("B13").Value = 'La maison en A'
$ws.Range("C13").Value = '2082'
$ws.Range("D13").Value = 'N/A'
$ws.Range("E13").Value = 'https://www.lego.com/cdn/cs/set/assets/blt2b163a472ef2e61f/21338.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1'
$ws.Range("F13").Value = 'https://www.lego.com/fr-fr/product/21338'
$ws.Range("G13").Value = ''
$ws.Range("H13").Value = ''
$ws.Range("I13").Value = ''
$ws.Range("J13").Value = ''

# --- Clear the Collection ("N/A") value for rows whose set now has no known collection ---
$clearDRows = 8,9,10,11,12,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28
foreach ($r in $clearDRows) {
    $ws.Range("D$r").Value = ''
}

$ws.Range("A1").Select()